# Generate Report for Handoff
#
# The "b.md" file has now been handed off again (new handback xliff files
# were generated), so its status flips from "Handed back: in sync with
# en-US" to "Ready for handoff" on every sheet, the "Content Duplicate"
# flag clears, the Latest Handoff File/Datetime move to the newly
# generated xliff, and an Error Detail note about the stale handback
# version is recorded.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e630a65d04678b2a30928ee0cbca9f2fcb76f2f8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5709807eb035e44678ec052d57090ee2cd654ab4/e2e/b.md."

# ---- Overview sheet: row 3 is the "b.md" row ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 20:37:26"

# ---- zh-cn sheet: row 3 is the "b.md" row ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text storage (matches source data, which stores
# "True"/"False" as plain text rather than real booleans).
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-27 20:37:20"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 is the "b.md" row ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-27 20:37:26"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
